$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-15 from 45233 to 45243
$ws.Range("C2:C15").Value = 45243
